# Generate Report for Handback
# Fills in the "Latest Target File" / "Latest Handback File" / "Latest Handback
# DateTime" columns on the per-language sheets now that the handback round
# trip has completed, flips the Overview/status columns from "Ready for
# handoff" to "Handed back: in sync with en-US", and widens the columns that
# now hold longer text.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

$statusText = "Handed back: in sync with en-US"

$baacdMdUrl  = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/c88c836f7dda94a04d244af0c335add5df96acbd/e2e/baacd752-d650-4f23-89e6-0db079eccae6.md"
$c1e58MdUrl  = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/c88c836f7dda94a04d244af0c335add5df96acbd/e2e/c1e58909-a01c-40a0-a0c9-266b2875f041.md"
$baacdMdName = "baacd752-d650-4f23-89e6-0db079eccae6.md"
$c1e58MdName = "c1e58909-a01c-40a0-a0c9-266b2875f041.md"

# ColumnWidth goes through Excel's pixel-quantised character-width engine
# before it lands in the saved XML `width`, so the inputs below are chosen
# (by solving that quantisation) to land as close as possible to the wider
# "29.9777047293527" / "40" columns used by the generated report.
$wideStatusColumnWidth = 29.166666666666668
$wideFileColumnWidth   = 39.166666666666664

# --- Overview sheet: zh-cn / de-de status columns (E, F) ---------------
$overview.Range("E2").Value = $statusText
$overview.Range("F2").Value = $statusText
$overview.Range("E3").Value = $statusText
$overview.Range("F3").Value = $statusText

$overview.Columns.Item(5).ColumnWidth = $wideStatusColumnWidth
$overview.Columns.Item(6).ColumnWidth = $wideStatusColumnWidth

# --- zh-cn sheet ---------------------------------------------------------
$zhcn.Range("C2").Value = $statusText
$zhcn.Range("C3").Value = $statusText
$zhcn.Columns.Item(3).ColumnWidth = $wideStatusColumnWidth

$zhcn.Range("K2").Value = "2016-08-17 16:46:50"
$zhcn.Range("K3").Value = "2016-08-17 16:46:50"

$zhcn.Range("J2").Value = "baacd752-d650-4f23-89e6-0db079eccae6.3eac469c0f8383a76f0040ae99bcc80501f82a8c.zh-cn.xlf"
$zhcn.Range("J3").Value = "c1e58909-a01c-40a0-a0c9-266b2875f041.e8562b93bdd49870a7773764ab055171c7c4c662.zh-cn.xlf"

$zhcn.Range("I2").Value = $baacdMdName
$zhcn.Range("I3").Value = $c1e58MdName

$zhcn.Columns.Item(9).ColumnWidth = $wideFileColumnWidth
$zhcn.Columns.Item(10).ColumnWidth = $wideFileColumnWidth

$zhcn.Hyperlinks.Delete()
$zhcn.Hyperlinks.Add($zhcn.Range("A2"), $baacdMdUrl, "", "", $baacdMdName)
$zhcn.Hyperlinks.Add($zhcn.Range("I2"), $baacdMdUrl, "", "", $baacdMdName)
$zhcn.Hyperlinks.Add($zhcn.Range("A3"), $c1e58MdUrl, "", "", $c1e58MdName)
$zhcn.Hyperlinks.Add($zhcn.Range("I3"), $c1e58MdUrl, "", "", $c1e58MdName)

# --- de-de sheet -----------------------------------------------------------
$dede.Range("C2").Value = $statusText
$dede.Range("C3").Value = $statusText
$dede.Columns.Item(3).ColumnWidth = $wideStatusColumnWidth

$dede.Range("K2").Value = "2016-08-17 16:46:58"
$dede.Range("K3").Value = "2016-08-17 16:46:58"

$dede.Range("J2").Value = "baacd752-d650-4f23-89e6-0db079eccae6.3eac469c0f8383a76f0040ae99bcc80501f82a8c.de-de.xlf"
$dede.Range("J3").Value = "c1e58909-a01c-40a0-a0c9-266b2875f041.e8562b93bdd49870a7773764ab055171c7c4c662.de-de.xlf"

$dede.Range("I2").Value = $baacdMdName
$dede.Range("I3").Value = $c1e58MdName

$dede.Columns.Item(9).ColumnWidth = $wideFileColumnWidth
$dede.Columns.Item(10).ColumnWidth = $wideFileColumnWidth

$dede.Hyperlinks.Delete()
$dede.Hyperlinks.Add($dede.Range("A2"), $baacdMdUrl, "", "", $baacdMdName)
$dede.Hyperlinks.Add($dede.Range("I2"), $baacdMdUrl, "", "", $baacdMdName)
$dede.Hyperlinks.Add($dede.Range("A3"), $c1e58MdUrl, "", "", $c1e58MdName)
$dede.Hyperlinks.Add($dede.Range("I3"), $c1e58MdUrl, "", "", $c1e58MdName)

Write-Host "Handback report generated"
